$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Formula = "=CONCAT(A4, B4)"
$ws.Range("D5").Formula = "=LEN(B5)"
$ws.Range("D6").Formula = "=LEFT(B6)"
$ws.Range("D7").Formula = "=RIGHT(B7)"
$ws.Range("D8").Formula = "=MID(B8, 1, 1)"
$ws.Range("D9").Formula = "=UPPER(B9)"
$ws.Range("D10").Formula = "=LOWER(B10)"
$ws.Range("D11").Formula = "=PROPER(B11)"
$ws.Range("D12").Formula = "=TRIM(B12)"
$ws.Range("D13").Formula = '=SUBSTITUTE(B13,"Iphones","Androids")'
$ws.Range("D14").Formula = '=REPLACE(B14, 1, 1, "You")'

$ws.Range("D19").Select()
